# [service] add appointment feature
#
# Renames the two FK columns on the "like" junction table to reflect that
# they are now also part of its composite primary key, and documents the
# new Flyway "schema_version" bookkeeping table (appointment feature
# migration tracking) in a new column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "like" table: user_a_id / user_b_id are now part of the PK as well ---
$ws.Range("D15").Value2 = "user_a_id (FK,PK)"
$ws.Range("D16").Value2 = "user_b_id (FK,PK)"

# --- New column F: Flyway schema_version table layout ---
$ws.Range("F1").Value2  = "schema_version"
$ws.Range("F2").Value2  = "installed_rank(PK)"
$ws.Range("F3").Value2  = "version"
$ws.Range("F4").Value2  = "description"
$ws.Range("F5").Value2  = "type"
$ws.Range("F6").Value2  = "script"
$ws.Range("F7").Value2  = "checksum"
$ws.Range("F8").Value2  = "installed_by"
$ws.Range("F9").Value2  = "installed_on"
$ws.Range("F10").Value2 = "execution_time"
$ws.Range("F11").Value2 = "success"

# Match the formatting of the other table headers/body cells in the sheet:
#  - F1 is the table title -> bold, black (same style as A1/B1/C1/D1)
#  - F2/F3 and F8/F9 are PK-ish header rows -> black, non-bold (same as A2/A3 etc.)
#  - F4-F7, F10, F11 stay on the default style (unchanged, like most body cells)
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Font.Color = 0

$ws.Range("F2").Font.Bold = $false
$ws.Range("F2").Font.Color = 0
$ws.Range("F3").Font.Bold = $false
$ws.Range("F3").Font.Color = 0

$ws.Range("F8").Font.Bold = $false
$ws.Range("F8").Font.Color = 0
$ws.Range("F9").Font.Bold = $false
$ws.Range("F9").Font.Color = 0

# --- Selection moved to B8 (matches the author's last editing position) ---
[void]$ws.Range("B8").Select()
